# Weekly update: insert a new latest-week record at row 33 (Alcachofa / Agrícola del Norte S.A. de Arica)
# This pushes the previous rows 33-46 down to 34-47 and fills the new row 33 with the
# newest observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 33; Excel shifts rows 33:46 down to 34:47
# and carries the formatting (including the date style on column D) down with them.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly record.
$ws.Range("A33").Value = 1
$ws.Range("B33").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C33").Value = "Arica y Parinacota"
$ws.Range("D33").Value = 45215
$ws.Range("E33").Value = 15
$ws.Range("F33").Value = 100112013
$ws.Range("G33").Value = "Alcachofa"
$ws.Range("H33").Value = "Madrigal"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 14000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 14500
$ws.Range("N33").Value = "$/caja 40 unidades"
$ws.Range("O33").Value = "Región de Coquimbo"
$ws.Range("P33").Value = 362
$ws.Range("Q33").Value = 40
$ws.Range("R33").Value = "Hortaliza"
